$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data gained two more reported quarters (Dec-2018 and Sep-2018).
# Insert two blank columns before column D; this shifts the existing
# D:K data (values, formats, styles) right to F:M automatically.
$ws.Range("D:E").EntireColumn.Insert()

# The freshly inserted D:E columns come back with the default "General"
# format, so re-apply the same number formats used by the rest of the
# table: date format for the "Period Ending" header rows, and the
# thousands-separator numeric format everywhere else.
$ws.Range("D7:E7,D38:E38,D80:E80").NumberFormat = "[$-409]d\-mmm\-yy;@"
$ws.Range("D8:E35,D39:E77,D81:E102").NumberFormat = "#,##0"

# Populate the two new quarters of data (columns D and E)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 57300
$ws.Range("E8").Value = 28100
$ws.Range("D9").Value = 38300
$ws.Range("E9").Value = 28300
$ws.Range("D10").Value = 19000
$ws.Range("E10").Value = -200
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 139900
$ws.Range("E17").Value = 122200
$ws.Range("D18").Value = -82600
$ws.Range("E18").Value = -94100
$ws.Range("D20").Value = 5000
$ws.Range("E20").Value = 4400
$ws.Range("D21").Value = -77600
$ws.Range("E21").Value = -89700
$ws.Range("D22").Value = 1600
$ws.Range("E22").Value = 900
$ws.Range("D23").Value = -79200
$ws.Range("E23").Value = -90600
$ws.Range("D24").Value = 1800
$ws.Range("E24").Value = -400
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -81000
$ws.Range("E26").Value = -90300
$ws.Range("D27").Value = -80800
$ws.Range("E27").Value = -89100
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -5000
$ws.Range("E32").Value = -4400
$ws.Range("D33").Value = -80800
$ws.Range("E33").Value = -89100
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -80800
$ws.Range("E35").Value = -89100
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 18200
$ws.Range("E41").Value = 43900
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 69200
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 283900
$ws.Range("E47").Value = 478000
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 685000
$ws.Range("E52").Value = 673800
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1435400
$ws.Range("E54").Value = 1704300
$ws.Range("D57").Value = 0
$ws.Range("E57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 351000
$ws.Range("E59").Value = 429300
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 91200
$ws.Range("E61").Value = 89600
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 958200
$ws.Range("E66").Value = 1147300
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -26100
$ws.Range("E72").Value = 54700
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 477300
$ws.Range("E76").Value = 557000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -80800
$ws.Range("E81").Value = -89100
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -37200
$ws.Range("E89").Value = 6600
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 23100
$ws.Range("E94").Value = -701300
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("E100").Value = 82800
$ws.Range("D101").Value = -400
$ws.Range("E101").Value = -700
$ws.Range("D102").Value = -14500
$ws.Range("E102").Value = -612600

# Data corrections beyond the simple column shift
# Row 22 (Interest Expense): F:J become "NA" (previously 0)
$ws.Range("F22:J22").Value = "NA"

# Row 94 (Total Cash Flows From Investing Activities): H/I corrected values
$ws.Range("H94").Value = 128700
$ws.Range("I94").Value = -40200

# Row 102 (Change In Cash and Cash Equivalents): H/I corrected values
$ws.Range("H102").Value = 149200
$ws.Range("I102").Value = -12600
